# Update "想去人数" (interested-people count) figures for several events.
# These values appear identically on both the "展览" sheet and the
# "全部类型" sheet, so both need to be updated.

$wb = $excel.ActiveWorkbook

# Row -> new value for column F
$updates = @{
    4  = 2879
    7  = 13
    8  = 1538
    12 = 1276
    14 = 409
    21 = 2869
    22 = 346
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
